$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# D-column cells are forced to Text format while assigning so that numeric-looking
# strings (e.g. "583.86") are stored as text, not auto-converted to numbers by Excel;
# the cell style is restored immediately afterward so no new style index is introduced.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.238.16"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -0.87%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.504.22"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.06%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.86"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.56%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.64"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.95%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.505.71"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.28%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  -2.41%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.105.79"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -0.07%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.29"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -1.71%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.79%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.502.62"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.36%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.271.34"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -0.84%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.79"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -2.05%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.88"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -2.75%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.59"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.42%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.85"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -2.13%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.50%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.645.62"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -0.21%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.88"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  -0.06%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.73"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +0.10%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000115"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("E31").Value = "  +0.10%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.33"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  -2.10%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.519.14"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +0.21%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.60"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.55%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.30"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  -2.64%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.88"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -1.47%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "164.12"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("E42").Value = "  -2.79%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.809"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -0.77%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.13"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  +0.12%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.78"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("E49").Value = "  -1.02%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.478.45"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.30%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.76"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -1.69%  "
